# Update the cached "Date automatically updates" text shown in the
# Date Placeholder shape that lives on the Slide Master and on every
# Slide Layout (Insert > Header & Footer > Date and time).
#
# The deck's slides themselves do not carry their own date placeholder
# shape (they inherit it from their layout), so touching the Master +
# every CustomLayout covers every occurrence of the cached date text.

$p = $ppt.ActivePresentation

$newDate = "8/26/14"

function Update-DatePlaceholders($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master's own Date Placeholder shape.
$master = $p.SlideMaster
Update-DatePlaceholders $master

# Every Slide Layout hanging off the master has its own Date
# Placeholder shape too.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholders $layout
}
